# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Mon Apr 22 15:21:45 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.781.27"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "3.173.71"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.170.08"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "3.693.64"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "65.872.79"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "3.181.13"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0893"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "479.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.01%  "
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0638"
$ws.Range("E44").Value = "  +8.12%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.878.62"
$ws.Range("E45").Value = "  -5.96%  "
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.31%  "
